$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 950

# Row 3
$ws.Range("H3").Value = 2.16
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 4.9
$ws.Range("R3").Value = 1.29
$ws.Range("S3").Value = 3.3

# Row 4
$ws.Range("F4").Value = 1.69
$ws.Range("H4").Value = 4.5
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 980
$ws.Range("N4").Value = 2.46
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 1.46

# Row 5
$ws.Range("F5").Value = 2.36
$ws.Range("G5").Value = 2.86
$ws.Range("H5").Value = 2.52
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3.95
$ws.Range("K5").Value = 5.3
$ws.Range("L5").Value = 1.2
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 2.8
$ws.Range("O5").Value = 1.13
$ws.Range("P5").Value = 2.8
$ws.Range("Q5").Value = 1.44
$ws.Range("R5").Value = 1.63
$ws.Range("S5").Value = 1.93
$ws.Range("T5").Value = 1.01
$ws.Range("U5").Value = 1.01
$ws.Range("V5").Value = 1.52
$ws.Range("W5").Value = 1.58
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# Row 6
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 2.54
$ws.Range("O6").Value = 1.17
$ws.Range("R6").Value = 1.2
$ws.Range("S6").Value = 2.04
$ws.Range("T6").Value = 1.01
$ws.Range("U6").Value = 1.01
$ws.Range("V6").Value = 1.04
$ws.Range("W6").Value = 4.6
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 1000
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000

# Row 7
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 1.89
$ws.Range("O7").Value = 1.01
$ws.Range("R7").Value = 1.1
$ws.Range("S7").Value = 1.89
$ws.Range("T7").Value = 1.01
$ws.Range("U7").Value = 1.01
$ws.Range("V7").Value = 3.1
$ws.Range("W7").Value = 1.06
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 1000
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 1000
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 1000

# Row 8
$ws.Range("H8").Value = 2.3
$ws.Range("L8").Value = 1.44
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 3.15
$ws.Range("O8").Value = 1.36
$ws.Range("R8").Value = 1.28
$ws.Range("S8").Value = 3.75
$ws.Range("T8").Value = 1.81
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.63
$ws.Range("W8").Value = 1.38
$ws.Range("X8").Value = 14.5
$ws.Range("Y8").Value = 11.5
$ws.Range("Z8").Value = 18
$ws.Range("AA8").Value = 42
$ws.Range("AB8").Value = 14.5
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 13.5
$ws.Range("AE8").Value = 36
$ws.Range("AF8").Value = 28
$ws.Range("AG8").Value = 17
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 55
$ws.Range("AJ8").Value = 75
$ws.Range("AK8").Value = 55
$ws.Range("AL8").Value = 60
$ws.Range("AM8").Value = 140
$ws.Range("AN8").Value = 55
$ws.Range("AO8").Value = 30

# Row 9
$ws.Range("G9").Value = 2.36
$ws.Range("I9").Value = 3.35
$ws.Range("J9").Value = 3.65
$ws.Range("L9").Value = 1.34
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 4.6
$ws.Range("O9").Value = 1.23
$ws.Range("P9").Value = 2.2
$ws.Range("Q9").Value = 1.68
$ws.Range("R9").Value = 1.46
$ws.Range("S9").Value = 2.6
$ws.Range("T9").Value = 1.57
$ws.Range("U9").Value = 2.38
$ws.Range("V9").Value = 1.42
$ws.Range("W9").Value = 1.73
$ws.Range("X9").Value = 19
$ws.Range("Y9").Value = 17
$ws.Range("Z9").Value = 26
$ws.Range("AA9").Value = 55
$ws.Range("AB9").Value = 13.5
$ws.Range("AC9").Value = 8.6
$ws.Range("AD9").Value = 15
$ws.Range("AE9").Value = 34
$ws.Range("AF9").Value = 17.5
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 15
$ws.Range("AI9").Value = 40
$ws.Range("AJ9").Value = 32
$ws.Range("AK9").Value = 23
$ws.Range("AL9").Value = 32
$ws.Range("AM9").Value = 980
$ws.Range("AN9").Value = 15
$ws.Range("AO9").Value = 27

# Row 10
$ws.Range("I10").Value = 5.5
$ws.Range("K10").Value = 3.75
$ws.Range("L10").Value = 1.42
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 3.1
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 1.86
$ws.Range("Q10").Value = 1.92
$ws.Range("R10").Value = 1.25
$ws.Range("S10").Value = 3.5
$ws.Range("T10").Value = 1.8
$ws.Range("U10").Value = 1.98
$ws.Range("V10").Value = 1.23
$ws.Range("W10").Value = 2.1
$ws.Range("X10").Value = 17
$ws.Range("Y10").Value = 17.5
$ws.Range("Z10").Value = 44
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 9
$ws.Range("AC10").Value = 8.800000000000001
$ws.Range("AD10").Value = 21
$ws.Range("AE10").Value = 70
$ws.Range("AF10").Value = 12
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 20
$ws.Range("AI10").Value = 75
$ws.Range("AJ10").Value = 22
$ws.Range("AK10").Value = 22
$ws.Range("AL10").Value = 40
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 14
$ws.Range("AO10").Value = 80

# Row 11
$ws.Range("F11").Value = 1.74
$ws.Range("G11").Value = 1.93
$ws.Range("H11").Value = 4.8
$ws.Range("I11").Value = 6
$ws.Range("J11").Value = 3.4
$ws.Range("K11").Value = 4.6
$ws.Range("L11").Value = 1.42
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 3.35
$ws.Range("O11").Value = 1.35
$ws.Range("P11").Value = 1.79
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = 1.3
$ws.Range("S11").Value = 3.3
$ws.Range("T11").Value = 1.77
$ws.Range("U11").Value = 1.76
$ws.Range("V11").Value = 1.2
$ws.Range("W11").Value = 2.06
$ws.Range("X11").Value = 16
$ws.Range("Y11").Value = 980
$ws.Range("Z11").Value = 980
$ws.Range("AA11").Value = 170
$ws.Range("AB11").Value = 9.6
$ws.Range("AC11").Value = 10.5
$ws.Range("AD11").Value = 980
$ws.Range("AE11").Value = 95
$ws.Range("AF11").Value = 13
$ws.Range("AG11").Value = 12.5
$ws.Range("AH11").Value = 980
$ws.Range("AI11").Value = 100
$ws.Range("AJ11").Value = 980
$ws.Range("AK11").Value = 980
$ws.Range("AL11").Value = 980
$ws.Range("AM11").Value = 160
$ws.Range("AN11").Value = 17
$ws.Range("AO11").Value = 130

# Row 12
$ws.Range("F12").Value = 2.94
$ws.Range("G12").Value = 3.8
$ws.Range("H12").Value = 2.2
$ws.Range("I12").Value = 2.66
$ws.Range("J12").Value = 3.55
$ws.Range("K12").Value = 4.2
$ws.Range("P12").Value = 2.16
$ws.Range("Q12").Value = 1.68

# Row 13
$ws.Range("F13").Value = 1.28
$ws.Range("G13").Value = 1.32
$ws.Range("H13").Value = 10.5
$ws.Range("J13").Value = 6.6
$ws.Range("K13").Value = 7

# Row 14
$ws.Range("F14").Value = 2.78
$ws.Range("G14").Value = 3.25
$ws.Range("H14").Value = 2.24
$ws.Range("I14").Value = 2.78
$ws.Range("J14").Value = 3.35
$ws.Range("K14").Value = 5.1
$ws.Range("P14").Value = 2.44
$ws.Range("Q14").Value = 1.55

# Row 15
$ws.Range("F15").Value = 3.05
$ws.Range("G15").Value = 3.85
$ws.Range("H15").Value = 1.97
$ws.Range("I15").Value = 2.28
$ws.Range("J15").Value = 4.3
$ws.Range("K15").Value = 5.9
$ws.Range("P15").Value = 3.55
$ws.Range("Q15").Value = 1.28

# Row 18
$ws.Range("H18").Value = 8.6
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 4.5

# Row 19
$ws.Range("I19").Value = 3.85
$ws.Range("R19").Value = 1.16
$ws.Range("S19").Value = 6
$ws.Range("T19").Value = 1.99
$ws.Range("U19").Value = 1.59
$ws.Range("V19").Value = 1.35
$ws.Range("X19").Value = 8
$ws.Range("Y19").Value = 9.800000000000001
$ws.Range("Z19").Value = 980
$ws.Range("AA19").Value = 100
$ws.Range("AB19").Value = 8.800000000000001
$ws.Range("AC19").Value = 7.4
$ws.Range("AD19").Value = 980
$ws.Range("AE19").Value = 70
$ws.Range("AF19").Value = 18
$ws.Range("AG19").Value = 14
$ws.Range("AH19").Value = 980
$ws.Range("AI19").Value = 120
$ws.Range("AJ19").Value = 980
$ws.Range("AK19").Value = 980
$ws.Range("AL19").Value = 80
$ws.Range("AM19").Value = 270
$ws.Range("AN19").Value = 980
$ws.Range("AO19").Value = 120

# Row 20
$ws.Range("G20").Value = 2.08
$ws.Range("H20").Value = 4.3
$ws.Range("AJ20").Value = 24

# Row 21
$ws.Range("F21").Value = 6.8
$ws.Range("G21").Value = 7.8
$ws.Range("I21").Value = 1.62
$ws.Range("J21").Value = 4.1
$ws.Range("K21").Value = 4.8
$ws.Range("P21").Value = 2
$ws.Range("Q21").Value = 1.82

# Row 22
$ws.Range("F22").Value = 1.97
$ws.Range("G22").Value = 2.14
$ws.Range("H22").Value = 3.8
$ws.Range("I22").Value = 4.5
$ws.Range("K22").Value = 4
$ws.Range("L22").Value = 1.33
$ws.Range("Q22").Value = 1.88
$ws.Range("T22").Value = 1.74
$ws.Range("V22").Value = 1.29
$ws.Range("W22").Value = 1.87
$ws.Range("Y22").Value = 18
$ws.Range("AB22").Value = 10.5
$ws.Range("AC22").Value = 9.6
$ws.Range("AD22").Value = 17.5
$ws.Range("AE22").Value = 60
$ws.Range("AF22").Value = 14
$ws.Range("AG22").Value = 11
$ws.Range("AI22").Value = 60
$ws.Range("AJ22").Value = 29
$ws.Range("AK22").Value = 23
$ws.Range("AL22").Value = 38
$ws.Range("AN22").Value = 16.5
